$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4571641683578491
$ws.Range("B1").Value = 0.9303494095802307
$ws.Range("C1").Value = 1.03879976272583
$ws.Range("D1").Value = 5.127120494842529
$ws.Range("E1").Value = 1.270225882530212
